$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete columns E:F (old "부서" and "직무" columns), shifting G,H left to E,F
$ws.Range("E:F").Delete()

# Rename remaining headers
$ws.Range("C4").Value = "최소자리"
$ws.Range("D4").Value = "복잡성"

# Update selection to match the recorded end state
$ws.Range("O27").Select()
